$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newDate = 42731.0647456713

$ws.Range("B3:G3").UnMerge()

$ws.Cells.Item(3, 2).Value = $newDate
$ws.Cells.Item(3, 3).Value = $newDate
$ws.Cells.Item(3, 4).Value = $newDate
$ws.Cells.Item(3, 5).Value = $newDate
$ws.Cells.Item(3, 6).Value = $newDate
$ws.Cells.Item(3, 7).Value = $newDate

$ws.Range("B3:G3").Merge()

$ws.Cells.Item(3, 2).Style = "FirstLineStyle"
